$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update time_taken (column F) timestamps on the "data" sheet
$ws1.Range("F2").Value = "2021-10-05 14:19:09.766768"
$ws1.Range("F3").Value = "2021-10-05 14:19:09.766776"
$ws1.Range("F4").Value = "2021-10-05 14:19:09.766779"
$ws1.Range("F5").Value = "2021-10-05 14:19:09.766782"
$ws1.Range("F6").Value = "2021-10-05 14:19:09.766785"
$ws1.Range("F7").Value = "2021-10-05 14:19:09.766788"
$ws1.Range("F8").Value = "2021-10-05 14:19:09.766790"
$ws1.Range("F9").Value = "2021-10-05 14:19:09.766792"
$ws1.Range("F10").Value = "2021-10-05 14:19:09.766795"
$ws1.Range("F11").Value = "2021-10-05 14:19:09.766798"
$ws1.Range("F12").Value = "2021-10-05 14:19:09.766800"
$ws1.Range("F13").Value = "2021-10-05 14:19:09.766803"
$ws1.Range("F14").Value = "2021-10-05 14:19:09.766805"
$ws1.Range("F15").Value = "2021-10-05 14:19:09.766808"
$ws1.Range("F16").Value = "2021-10-05 14:19:09.766810"
$ws1.Range("F17").Value = "2021-10-05 14:19:09.766813"
$ws1.Range("F18").Value = "2021-10-05 14:19:09.766815"
$ws1.Range("F19").Value = "2021-10-05 14:19:09.766818"
$ws1.Range("F20").Value = "2021-10-05 14:19:09.766820"
$ws1.Range("F21").Value = "2021-10-05 14:19:09.766823"
$ws1.Range("F22").Value = "2021-10-05 14:19:09.766825"
$ws1.Range("F23").Value = "2021-10-05 14:19:09.766827"
$ws1.Range("F24").Value = "2021-10-05 14:19:09.766830"
$ws1.Range("F25").Value = "2021-10-05 14:19:09.766832"
$ws1.Range("F26").Value = "2021-10-05 14:19:09.766835"
$ws1.Range("F27").Value = "2021-10-05 14:19:09.766838"
$ws1.Range("F28").Value = "2021-10-05 14:19:09.766840"
$ws1.Range("F29").Value = "2021-10-05 14:19:09.766843"
$ws1.Range("F30").Value = "2021-10-05 14:19:09.766845"
$ws1.Range("F31").Value = "2021-10-05 14:19:09.766847"
$ws1.Range("F32").Value = "2021-10-05 14:19:09.766850"
$ws1.Range("F33").Value = "2021-10-05 14:19:09.766852"
$ws1.Range("F34").Value = "2021-10-05 14:19:09.766855"
$ws1.Range("F35").Value = "2021-10-05 14:19:09.766858"
$ws1.Range("F36").Value = "2021-10-05 14:19:09.766860"
$ws1.Range("F37").Value = "2021-10-05 14:19:09.766863"
$ws1.Range("F38").Value = "2021-10-05 14:19:09.766865"
$ws1.Range("F39").Value = "2021-10-05 14:19:09.766868"
$ws1.Range("F40").Value = "2021-10-05 14:19:09.766870"
$ws1.Range("F41").Value = "2021-10-05 14:19:09.766872"
$ws1.Range("F42").Value = "2021-10-05 14:19:09.766875"
$ws1.Range("F43").Value = "2021-10-05 14:19:09.766878"
$ws1.Range("F44").Value = "2021-10-05 14:19:09.766880"
$ws1.Range("F45").Value = "2021-10-05 14:19:09.766883"
$ws1.Range("F46").Value = "2021-10-05 14:19:09.766885"
$ws1.Range("F47").Value = "2021-10-05 14:19:09.766888"
$ws1.Range("F48").Value = "2021-10-05 14:19:09.766890"
$ws1.Range("F49").Value = "2021-10-05 14:19:09.766892"
$ws1.Range("F50").Value = "2021-10-05 14:19:09.766895"
$ws1.Range("F51").Value = "2021-10-05 14:19:09.766897"
$ws1.Range("F52").Value = "2021-10-05 14:19:09.766900"
$ws1.Range("F53").Value = "2021-10-05 14:19:09.766902"
$ws1.Range("F54").Value = "2021-10-05 14:19:09.766905"
$ws1.Range("F55").Value = "2021-10-05 14:19:09.766907"
$ws1.Range("F56").Value = "2021-10-05 14:19:09.766910"
$ws1.Range("F57").Value = "2021-10-05 14:19:09.766912"
$ws1.Range("F58").Value = "2021-10-05 14:19:09.766915"
$ws1.Range("F59").Value = "2021-10-05 14:19:09.766917"
$ws1.Range("F60").Value = "2021-10-05 14:19:09.766920"
$ws1.Range("F61").Value = "2021-10-05 14:19:09.766922"
$ws1.Range("F62").Value = "2021-10-05 14:19:09.766925"
$ws1.Range("F63").Value = "2021-10-05 14:19:09.766927"
$ws1.Range("F64").Value = "2021-10-05 14:19:09.766930"

# Add a new "metadata" worksheet positioned after "data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Match sheetPr/outline + page margin settings used elsewhere in this workbook
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Copy header cell formatting (bold + border + center/top alignment) from the
# "data" sheet's header row so the new sheet matches the workbook's existing style.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Anophthalmia or microphthalmia"
$ws2.Range("C2").Value = 34

# data_version must stay a text value ("1.42") rather than become a number
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.42"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("E2").Value = "2021-09-09T11:07:31.123840Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:09.763112"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/34/?format=json"

$ws1.Select()
